# Removed cc to support email
# Append " *" to each header cell (A1:E1) on the active sheet, then leave
# the selection on A2 (matching the recorded end-user selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:E1")
foreach ($cell in $headerRange.Cells) {
    $cell.Value = "$($cell.Value2) *"
}

$ws.Range("A2").Select()
